$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, shifting existing rows 38-78 down to 39-79.
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the new record's data.
$ws.Range("A38").Value2 = 5
$ws.Range("B38").Value2 = "Macroferia Regional de Talca"
$ws.Range("C38").Value2 = "Maule"
$ws.Range("D38").Value2 = 44539
$ws.Range("E38").Value2 = 7
$ws.Range("F38").Value2 = 100112022
$ws.Range("G38").Value2 = "Arveja Verde"
$ws.Range("H38").Value2 = "Sin especificar"
$ws.Range("I38").Value2 = "Primera"
$ws.Range("J38").Value2 = 200
$ws.Range("K38").Value2 = 18000
$ws.Range("L38").Value2 = 18000
$ws.Range("M38").Value2 = 18000
$ws.Range("N38").Value2 = "`$/saco 25 kilos"
$ws.Range("O38").Value2 = "Región de La Araucanía"
$ws.Range("P38").Value2 = 720
$ws.Range("Q38").Value2 = 25
$ws.Range("R38").Value2 = "Hortaliza"
